# "Changed the emergence number" - update the EMER FUND budget line (B34)
# in the AMFOSTER BUDGET 2024 workbook. Dependent totals (B36 = SUM(B33:B35)
# and B39 = B7-B28-B36) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("B34").Value = 500

# Reflect the saved view state: scroll the window and move the selection.
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A15").Select()
